$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40 (Leve Item ID 5505)
$ws.Cells.Item(40, 8).Value = 1500.2  # H40: 1716.3636 -> 1500.2
$ws.Cells.Item(40, 9).Value = 1316.4615  # I40: 1481.6666 -> 1316.4615
$ws.Cells.Item(40, 10).Value = 1841.4286  # J40: 1998 -> 1841.4286
$ws.Cells.Item(40, 11).Value = 1316.4615  # K40: 1481.6666 -> 1316.4615
$ws.Cells.Item(40, 12).Value = 1841.4286  # L40: 1998 -> 1841.4286
$ws.Cells.Item(40, 13).Value = -1141.4615  # M40: -1306.6666 -> -1141.4615
$ws.Cells.Item(40, 14).Value = -2191.4286  # N40: -2348 -> -2191.4286
# Row 114 (Leve Item ID 25959)
$ws.Cells.Item(114, 8).Value = 0  # H114: 36000 -> 0
$ws.Cells.Item(114, 10).Value = 0  # J114: 36000 -> 0
$ws.Cells.Item(114, 12).ClearContents()  # L114: 36000 -> (removed)
$ws.Cells.Item(114, 14).Value = 0  # N114: -44678 -> 0
# Row 129 (Leve Item ID 36115)
$ws.Cells.Item(129, 8).Value = 897.48  # H129: 898.1799999999999 -> 897.48
$ws.Cells.Item(129, 10).Value = 897.48  # J129: 898.1799999999999 -> 897.48
$ws.Cells.Item(129, 12).Value = 2692.44  # L129: 2694.54 -> 2692.44
$ws.Cells.Item(129, 14).Value = -12692.44  # N129: -12694.54 -> -12692.44

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Cells.Item(2, 8).Value = 975.44116  # H2: 1168.5238 -> 975.44116
$ws.Cells.Item(2, 9).Value = 963.6774  # I2: 1161.95 -> 963.6774
$ws.Cells.Item(2, 10).Value = 1097  # J2: 1300 -> 1097
$ws.Cells.Item(2, 11).Value = 963.6774  # K2: 1161.95 -> 963.6774
$ws.Cells.Item(2, 12).Value = 1097  # L2: 1300 -> 1097
$ws.Cells.Item(2, 13).Value = -850.6774  # M2: -1048.95 -> -850.6774
$ws.Cells.Item(2, 14).Value = -1323  # N2: -1526 -> -1323
# Row 32 (Leve Item ID 44147)
$ws.Cells.Item(32, 8).Value = 9260.035  # H32: 9721.305 -> 9260.035
$ws.Cells.Item(32, 9).Value = 9889.067999999999  # I32: 10601.412 -> 9889.067999999999
$ws.Cells.Item(32, 10).Value = 5727.769  # J32: 5446.5 -> 5727.769
$ws.Cells.Item(32, 11).Value = 9889.067999999999  # K32: 10601.412 -> 9889.067999999999
$ws.Cells.Item(32, 12).Value = 5727.769  # L32: 5446.5 -> 5727.769
$ws.Cells.Item(32, 13).Value = -9602.067999999999  # M32: -10314.412 -> -9602.067999999999
$ws.Cells.Item(32, 14).Value = -6301.769  # N32: -6020.5 -> -6301.769
# Row 61 (Leve Item ID 43999)
$ws.Cells.Item(61, 8).Value = 3023.3635  # H61: 2147.2632 -> 3023.3635
$ws.Cells.Item(61, 9).Value = 3231.8462  # I61: 2040.64 -> 3231.8462
$ws.Cells.Item(61, 10).Value = 2722.2222  # J61: 2352.3076 -> 2722.2222
$ws.Cells.Item(61, 11).Value = 3231.8462  # K61: 2040.64 -> 3231.8462
$ws.Cells.Item(61, 12).Value = 2722.2222  # L61: 2352.3076 -> 2722.2222
$ws.Cells.Item(61, 13).Value = -3019.8462  # M61: -1828.64 -> -3019.8462
$ws.Cells.Item(61, 14).Value = -3146.2222  # N61: -2776.3076 -> -3146.2222
# Row 101 (Leve Item ID 18518)
$ws.Cells.Item(101, 8).Value = 42933.332  # H101: 56602 -> 42933.332
$ws.Cells.Item(101, 10).Value = 42933.332  # J101: 56602 -> 42933.332
$ws.Cells.Item(101, 12).Value = 42933.332  # L101: 56602 -> 42933.332
$ws.Cells.Item(101, 14).Value = -49423.332  # N101: -63092 -> -49423.332
# Row 107 (Leve Item ID 25645)
$ws.Cells.Item(107, 8).Value = 0  # H107: 40000 -> 0
$ws.Cells.Item(107, 10).Value = 0  # J107: 40000 -> 0
$ws.Cells.Item(107, 12).ClearContents()  # L107: 40000 -> (removed)
$ws.Cells.Item(107, 14).Value = 0  # N107: -47680 -> 0
# Row 113 (Leve Item ID 26002)
$ws.Cells.Item(113, 8).Value = 39699  # H113: 40000 -> 39699
$ws.Cells.Item(113, 10).Value = 39699  # J113: 40000 -> 39699
$ws.Cells.Item(113, 12).Value = 39699  # L113: 40000 -> 39699
$ws.Cells.Item(113, 14).Value = -48377  # N113: -48678 -> -48377
# Row 114 (Leve Item ID 25968)
$ws.Cells.Item(114, 8).Value = 42399.332  # H114: 34099.5 -> 42399.332
$ws.Cells.Item(114, 10).Value = 42399.332  # J114: 34099.5 -> 42399.332
$ws.Cells.Item(114, 12).Value = 42399.332  # L114: 34099.5 -> 42399.332
$ws.Cells.Item(114, 14).Value = -51077.332  # N114: -42777.5 -> -51077.332
# Row 116 (Leve Item ID 27713)
$ws.Cells.Item(116, 8).Value = 975.44116  # H116: 1168.5238 -> 975.44116
$ws.Cells.Item(116, 9).Value = 963.6774  # I116: 1161.95 -> 963.6774
$ws.Cells.Item(116, 10).Value = 1097  # J116: 1300 -> 1097
$ws.Cells.Item(116, 11).Value = 963.6774  # K116: 1161.95 -> 963.6774
$ws.Cells.Item(116, 12).Value = 1097  # L116: 1300 -> 1097
$ws.Cells.Item(116, 13).Value = 1330.3226  # M116: 1132.05 -> 1330.3226
$ws.Cells.Item(116, 14).Value = -5685  # N116: -5888 -> -5685
# Row 119 (Leve Item ID 26287)
$ws.Cells.Item(119, 8).Value = 30590.455  # H119: 29856.715 -> 30590.455
$ws.Cells.Item(119, 10).Value = 30590.455  # J119: 29856.715 -> 30590.455
$ws.Cells.Item(119, 12).Value = 30590.455  # L119: 29856.715 -> 30590.455
$ws.Cells.Item(119, 14).Value = -40266.455  # N119: -39532.715 -> -40266.455
# Row 122 (Leve Item ID 36168)
$ws.Cells.Item(122, 8).Value = 6279.625  # H122: 6511.2256 -> 6279.625
$ws.Cells.Item(122, 9).Value = 6326.2144  # I122: 6593.852 -> 6326.2144
$ws.Cells.Item(122, 11).Value = 18978.6432  # K122: 19781.556 -> 18978.6432
$ws.Cells.Item(122, 13).Value = -16528.6432  # M122: -17331.556 -> -16528.6432
# Row 123 (Leve Item ID 34107)
$ws.Cells.Item(123, 8).Value = 25427.285  # H123: 25354.666 -> 25427.285
$ws.Cells.Item(123, 10).Value = 25427.285  # J123: 25354.666 -> 25427.285
$ws.Cells.Item(123, 12).Value = 25427.285  # L123: 25354.666 -> 25427.285
$ws.Cells.Item(123, 14).Value = -35227.285  # N123: -35154.666 -> -35227.285
# Row 125 (Leve Item ID 34251)
$ws.Cells.Item(125, 8).Value = 61428.75  # H125: 47857.5 -> 61428.75
$ws.Cells.Item(125, 10).Value = 61428.75  # J125: 47857.5 -> 61428.75
$ws.Cells.Item(125, 12).Value = 61428.75  # L125: 47857.5 -> 61428.75
$ws.Cells.Item(125, 14).Value = -71268.75  # N125: -57697.5 -> -71268.75
# Row 132 (Leve Item ID 43997)
$ws.Cells.Item(132, 8).Value = 8291.714  # H132: 8509.484 -> 8291.714
$ws.Cells.Item(132, 9).Value = 8209.157999999999  # I132: 8598.611000000001 -> 8209.157999999999
$ws.Cells.Item(132, 10).Value = 8389.75  # J132: 8402.532999999999 -> 8389.75
$ws.Cells.Item(132, 11).Value = 24627.474  # K132: 25795.833 -> 24627.474
$ws.Cells.Item(132, 12).Value = 25169.25  # L132: 25207.599 -> 25169.25
$ws.Cells.Item(132, 13).Value = -22097.474  # M132: -23265.833 -> -22097.474
$ws.Cells.Item(132, 14).Value = -30229.25  # N132: -30267.599 -> -30229.25
# Row 133 (Leve Item ID 41857)
$ws.Cells.Item(133, 8).Value = 34740  # H133: 41925.715 -> 34740
$ws.Cells.Item(133, 10).Value = 34740  # J133: 41925.715 -> 34740
$ws.Cells.Item(133, 12).Value = 34740  # L133: 41925.715 -> 34740
$ws.Cells.Item(133, 14).Value = -39800  # N133: -46985.715 -> -39800
# Row 136 (Leve Item ID 43999)
$ws.Cells.Item(136, 8).Value = 3023.3635  # H136: 2147.2632 -> 3023.3635
$ws.Cells.Item(136, 9).Value = 3231.8462  # I136: 2040.64 -> 3231.8462
$ws.Cells.Item(136, 10).Value = 2722.2222  # J136: 2352.3076 -> 2722.2222
$ws.Cells.Item(136, 11).Value = 9695.5386  # K136: 6121.92 -> 9695.5386
$ws.Cells.Item(136, 12).Value = 8166.6666  # L136: 7056.9228 -> 8166.6666
$ws.Cells.Item(136, 13).Value = -7145.5386  # M136: -3571.92 -> -7145.5386
$ws.Cells.Item(136, 14).Value = -13266.6666  # N136: -12156.9228 -> -13266.6666

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Cells.Item(3, 8).Value = 975.44116  # H3: 1168.5238 -> 975.44116
$ws.Cells.Item(3, 9).Value = 963.6774  # I3: 1161.95 -> 963.6774
$ws.Cells.Item(3, 10).Value = 1097  # J3: 1300 -> 1097
$ws.Cells.Item(3, 11).Value = 963.6774  # K3: 1161.95 -> 963.6774
$ws.Cells.Item(3, 12).Value = 1097  # L3: 1300 -> 1097
$ws.Cells.Item(3, 13).Value = -849.6774  # M3: -1047.95 -> -849.6774
$ws.Cells.Item(3, 14).Value = -1325  # N3: -1528 -> -1325
# Row 20 (Leve Item ID 14149)
$ws.Cells.Item(20, 8).Value = 202355.2  # H20: 63713.375 -> 202355.2
$ws.Cells.Item(20, 9).Value = 500900  # I20: 143535.42 -> 500900
$ws.Cells.Item(20, 10).Value = 3325.3333  # J20: 1629.5555 -> 3325.3333
$ws.Cells.Item(20, 11).Value = 500900  # K20: 143535.42 -> 500900
$ws.Cells.Item(20, 12).Value = 3325.3333  # L20: 1629.5555 -> 3325.3333
$ws.Cells.Item(20, 13).Value = -500653  # M20: -143288.42 -> -500653
$ws.Cells.Item(20, 14).Value = -3819.3333  # N20: -2123.5555 -> -3819.3333
# Row 22 (Leve Item ID 5092)
$ws.Cells.Item(22, 8).Value = 8486.666999999999  # H22: 7290 -> 8486.666999999999
$ws.Cells.Item(22, 9).Value = 10084  # I22: 10086 -> 10084
$ws.Cells.Item(22, 10).Value = 500  # J22: 300 -> 500
$ws.Cells.Item(22, 11).Value = 10084  # K22: 10086 -> 10084
$ws.Cells.Item(22, 12).Value = 500  # L22: 300 -> 500
$ws.Cells.Item(22, 13).Value = -9911  # M22: -9913 -> -9911
$ws.Cells.Item(22, 14).Value = -846  # N22: -646 -> -846
# Row 108 (Leve Item ID 25643)
$ws.Cells.Item(108, 8).Value = 27513.572  # H108: 30149.5 -> 27513.572
$ws.Cells.Item(108, 10).Value = 27513.572  # J108: 30149.5 -> 27513.572
$ws.Cells.Item(108, 12).Value = 27513.572  # L108: 30149.5 -> 27513.572
$ws.Cells.Item(108, 14).Value = -35193.572  # N108: -37829.5 -> -35193.572

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 19 (Leve Item ID 2233)
$ws.Cells.Item(19, 8).Value = 1050  # H19: 939 -> 1050
$ws.Cells.Item(19, 9).Value = 66.666664  # I19: 173.75 -> 66.666664
$ws.Cells.Item(19, 11).Value = 66.666664  # K19: 173.75 -> 66.666664
$ws.Cells.Item(19, 13).Value = 103.333336  # M19: -3.75 -> 103.333336
# Row 24 (Leve Item ID 2233)
$ws.Cells.Item(24, 8).Value = 1050  # H24: 939 -> 1050
$ws.Cells.Item(24, 9).Value = 66.666664  # I24: 173.75 -> 66.666664
$ws.Cells.Item(24, 11).Value = 66.666664  # K24: 173.75 -> 66.666664
$ws.Cells.Item(24, 13).Value = 103.333336  # M24: -3.75 -> 103.333336
# Row 31 (Leve Item ID 44023)
$ws.Cells.Item(31, 8).Value = 2174.5103  # H31: 2407.8838 -> 2174.5103
$ws.Cells.Item(31, 9).Value = 2262.963  # I31: 2690.3635 -> 2262.963
$ws.Cells.Item(31, 10).Value = 2065.9546  # J31: 2111.9524 -> 2065.9546
$ws.Cells.Item(31, 11).Value = 2262.963  # K31: 2690.3635 -> 2262.963
$ws.Cells.Item(31, 12).Value = 2065.9546  # L31: 2111.9524 -> 2065.9546
$ws.Cells.Item(31, 13).Value = -1967.963  # M31: -2395.3635 -> -1967.963
$ws.Cells.Item(31, 14).Value = -2655.9546  # N31: -2701.9524 -> -2655.9546
# Row 34 (Leve Item ID 44023)
$ws.Cells.Item(34, 8).Value = 2174.5103  # H34: 2407.8838 -> 2174.5103
$ws.Cells.Item(34, 9).Value = 2262.963  # I34: 2690.3635 -> 2262.963
$ws.Cells.Item(34, 10).Value = 2065.9546  # J34: 2111.9524 -> 2065.9546
$ws.Cells.Item(34, 11).Value = 2262.963  # K34: 2690.3635 -> 2262.963
$ws.Cells.Item(34, 12).Value = 2065.9546  # L34: 2111.9524 -> 2065.9546
$ws.Cells.Item(34, 13).Value = -2060.963  # M34: -2488.3635 -> -2060.963
$ws.Cells.Item(34, 14).Value = -2469.9546  # N34: -2515.9524 -> -2469.9546
# Row 51 (Leve Item ID 2039)
$ws.Cells.Item(51, 8).Value = 20000  # H51: 20621.285 -> 20000
$ws.Cells.Item(51, 10).Value = 0  # J51: 20724.834 -> 0
$ws.Cells.Item(51, 12).Value = 0  # L51: 20724.834 -> 0
$ws.Cells.Item(51, 14).ClearContents()  # N51: -22196.834 -> (removed)
# Row 61 (Leve Item ID 2039)
$ws.Cells.Item(61, 8).Value = 20000  # H61: 20621.285 -> 20000
$ws.Cells.Item(61, 10).Value = 0  # J61: 20724.834 -> 0
$ws.Cells.Item(61, 12).Value = 0  # L61: 20724.834 -> 0
$ws.Cells.Item(61, 14).ClearContents()  # N61: -21420.834 -> (removed)
# Row 99 (Leve Item ID 36198)
$ws.Cells.Item(99, 8).Value = 1667.35  # H99: 1716.909 -> 1667.35
$ws.Cells.Item(99, 9).Value = 1661.0625  # I99: 1717.3334 -> 1661.0625
$ws.Cells.Item(99, 10).Value = 1692.5  # J99: 1715 -> 1692.5
$ws.Cells.Item(99, 11).Value = 1661.0625  # K99: 1717.3334 -> 1661.0625
$ws.Cells.Item(99, 12).Value = 1692.5  # L99: 1715 -> 1692.5
$ws.Cells.Item(99, 13).Value = -163.0625  # M99: -219.3334 -> -163.0625
$ws.Cells.Item(99, 14).Value = -4688.5  # N99: -4711 -> -4688.5
# Row 126 (Leve Item ID 36198)
$ws.Cells.Item(126, 8).Value = 1667.35  # H126: 1716.909 -> 1667.35
$ws.Cells.Item(126, 9).Value = 1661.0625  # I126: 1717.3334 -> 1661.0625
$ws.Cells.Item(126, 10).Value = 1692.5  # J126: 1715 -> 1692.5
$ws.Cells.Item(126, 11).Value = 4983.1875  # K126: 5152.0002 -> 4983.1875
$ws.Cells.Item(126, 12).Value = 5077.5  # L126: 5145 -> 5077.5
$ws.Cells.Item(126, 13).Value = -2513.1875  # M126: -2682.0002 -> -2513.1875
$ws.Cells.Item(126, 14).Value = -10017.5  # N126: -10085 -> -10017.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 37 (Leve Item ID 9516)
$ws.Cells.Item(37, 8).Value = 54500  # H37: 59500 -> 54500
$ws.Cells.Item(37, 10).Value = 54500  # J37: 59500 -> 54500
$ws.Cells.Item(37, 12).Value = 163500  # L37: 178500 -> 163500
$ws.Cells.Item(37, 14).Value = -163724  # N37: -178724 -> -163724
# Row 68 (Leve Item ID 12895)
$ws.Cells.Item(68, 8).Value = 186182.4  # H68: 209239.38 -> 186182.4
$ws.Cells.Item(68, 10).Value = 1860  # J68: 1948.8889 -> 1860
$ws.Cells.Item(68, 12).Value = 5580  # L68: 5846.6667 -> 5580
$ws.Cells.Item(68, 14).Value = -7202  # N68: -7468.6667 -> -7202
# Row 71 (Leve Item ID 12895)
$ws.Cells.Item(71, 8).Value = 186182.4  # H71: 209239.38 -> 186182.4
$ws.Cells.Item(71, 10).Value = 1860  # J71: 1948.8889 -> 1860
$ws.Cells.Item(71, 12).Value = 16740  # L71: 17540.0001 -> 16740
$ws.Cells.Item(71, 14).Value = -24852  # N71: -25652.0001 -> -24852
# Row 107 (Leve Item ID 27838)
$ws.Cells.Item(107, 8).Value = 590.13  # H107: 568.66 -> 590.13
$ws.Cells.Item(107, 9).Value = 400  # I107: 367.8387 -> 400
$ws.Cells.Item(107, 10).Value = 656.93243  # J107: 658.8840300000001 -> 656.93243
$ws.Cells.Item(107, 11).Value = 1200  # K107: 1103.5161 -> 1200
$ws.Cells.Item(107, 12).Value = 1970.79729  # L107: 1976.65209 -> 1970.79729
$ws.Cells.Item(107, 13).Value = 720  # M107: 816.4838999999999 -> 720
$ws.Cells.Item(107, 14).Value = -5810.79729  # N107: -5816.65209 -> -5810.79729
# Row 113 (Leve Item ID 27843)
$ws.Cells.Item(113, 8).Value = 164606.31  # H113: 185772.36 -> 164606.31
$ws.Cells.Item(113, 9).Value = 263730.6  # I113: 227817.4 -> 263730.6
$ws.Cells.Item(113, 10).Value = 835.73914  # J113: 774.1 -> 835.73914
$ws.Cells.Item(113, 11).Value = 791191.7999999999  # K113: 683452.2 -> 791191.7999999999
$ws.Cells.Item(113, 12).Value = 2507.21742  # L113: 2322.3 -> 2507.21742
$ws.Cells.Item(113, 13).Value = -789021.7999999999  # M113: -681282.2 -> -789021.7999999999
$ws.Cells.Item(113, 14).Value = -6847.21742  # N113: -6662.3 -> -6847.21742
# Row 131 (Leve Item ID 36060)
$ws.Cells.Item(131, 8).Value = 14495321  # H131: 17546820 -> 14495321
$ws.Cells.Item(131, 9).Value = 990  # I131: 743.3333 -> 990
$ws.Cells.Item(131, 10).Value = 14927987  # J131: 18521602 -> 14927987
$ws.Cells.Item(131, 11).Value = 2970  # K131: 2229.9999 -> 2970
$ws.Cells.Item(131, 12).Value = 44783961  # L131: 55564806 -> 44783961
$ws.Cells.Item(131, 13).Value = 2070  # M131: 2810.0001 -> 2070
$ws.Cells.Item(131, 14).Value = -44794041  # N131: -55574886 -> -44794041

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Cells.Item(70, 8).Value = 534806.9  # H70: 292648.56 -> 534806.9
$ws.Cells.Item(70, 9).Value = 633770.7  # I70: 424054.4 -> 633770.7
$ws.Cells.Item(70, 10).Value = 7000  # J70: 5944.909 -> 7000
$ws.Cells.Item(70, 11).Value = 633770.7  # K70: 424054.4 -> 633770.7
$ws.Cells.Item(70, 12).Value = 7000  # L70: 5944.909 -> 7000
$ws.Cells.Item(70, 13).Value = -633500.7  # M70: -423784.4 -> -633500.7
$ws.Cells.Item(70, 14).Value = -7540  # N70: -6484.909 -> -7540
# Row 73 (Leve Item ID 14146)
$ws.Cells.Item(73, 8).Value = 534806.9  # H73: 292648.56 -> 534806.9
$ws.Cells.Item(73, 9).Value = 633770.7  # I73: 424054.4 -> 633770.7
$ws.Cells.Item(73, 10).Value = 7000  # J73: 5944.909 -> 7000
$ws.Cells.Item(73, 11).Value = 633770.7  # K73: 424054.4 -> 633770.7
$ws.Cells.Item(73, 12).Value = 7000  # L73: 5944.909 -> 7000
$ws.Cells.Item(73, 13).Value = -632834.7  # M73: -423118.4 -> -632834.7
$ws.Cells.Item(73, 14).Value = -8872  # N73: -7816.909 -> -8872
# Row 132 (Leve Item ID 44008)
$ws.Cells.Item(132, 8).Value = 4038.7273  # H132: 2692.4119 -> 4038.7273
$ws.Cells.Item(132, 9).Value = 3400  # I132: 2349.818 -> 3400
$ws.Cells.Item(132, 10).Value = 4805.2  # J132: 3320.5 -> 4805.2
$ws.Cells.Item(132, 11).Value = 10200  # K132: 7049.454000000001 -> 10200
$ws.Cells.Item(132, 12).Value = 14415.6  # L132: 9961.5 -> 14415.6
$ws.Cells.Item(132, 13).Value = -7670  # M132: -4519.454000000001 -> -7670
$ws.Cells.Item(132, 14).Value = -19475.6  # N132: -15021.5 -> -19475.6
# Row 135 (Leve Item ID 42006)
$ws.Cells.Item(135, 8).Value = 0  # H135: 48112.855 -> 0
$ws.Cells.Item(135, 10).Value = 0  # J135: 48112.855 -> 0
$ws.Cells.Item(135, 12).ClearContents()  # L135: 48112.855 -> (removed)
$ws.Cells.Item(135, 14).Value = 0  # N135: -58252.855 -> 0
# Row 137 (Leve Item ID 43226)
$ws.Cells.Item(137, 8).Value = 73452.73  # H137: 75858.336 -> 73452.73
$ws.Cells.Item(137, 10).Value = 73452.73  # J137: 75858.336 -> 73452.73
$ws.Cells.Item(137, 12).Value = 73452.73  # L137: 75858.336 -> 73452.73
$ws.Cells.Item(137, 14).Value = -83652.73  # N137: -86058.336 -> -83652.73

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 14 (Leve Item ID 3771)
$ws.Cells.Item(14, 8).Value = 40000  # H14: 50000 -> 40000
$ws.Cells.Item(14, 10).Value = 30000  # J14: 0 -> 30000
$ws.Cells.Item(14, 12).Value = 30000  # L14: 0 -> 30000
$ws.Cells.Item(14, 14).Value = -30344  # N14: None -> -30344
# Row 16 (Leve Item ID 5289)
$ws.Cells.Item(16, 8).Value = 3112.625  # H16: 2240.75 -> 3112.625
$ws.Cells.Item(16, 9).Value = 850.25  # I16: 746.8 -> 850.25
$ws.Cells.Item(16, 10).Value = 5375  # J16: 3307.8572 -> 5375
$ws.Cells.Item(16, 11).Value = 850.25  # K16: 746.8 -> 850.25
$ws.Cells.Item(16, 12).Value = 5375  # L16: 3307.8572 -> 5375
$ws.Cells.Item(16, 13).Value = -680.25  # M16: -576.8 -> -680.25
$ws.Cells.Item(16, 14).Value = -5715  # N16: -3647.8572 -> -5715
# Row 104 (Leve Item ID 18675)
$ws.Cells.Item(104, 8).Value = 21948.572  # H104: 20567.5 -> 21948.572
$ws.Cells.Item(104, 10).Value = 21948.572  # J104: 20567.5 -> 21948.572
$ws.Cells.Item(104, 12).Value = 21948.572  # L104: 20567.5 -> 21948.572
$ws.Cells.Item(104, 14).Value = -28936.572  # N104: -27555.5 -> -28936.572
# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 4332.122  # H132: 4368.162 -> 4332.122
$ws.Cells.Item(132, 9).Value = 4429.853  # I132: 4503.5454 -> 4429.853
$ws.Cells.Item(132, 10).Value = 3857.4285  # J132: 3251.25 -> 3857.4285
$ws.Cells.Item(132, 11).Value = 13289.559  # K132: 13510.6362 -> 13289.559
$ws.Cells.Item(132, 12).Value = 11572.2855  # L132: 9753.75 -> 11572.2855
$ws.Cells.Item(132, 13).Value = -10759.559  # M132: -10980.6362 -> -10759.559
$ws.Cells.Item(132, 14).Value = -16632.2855  # N132: -14813.75 -> -16632.2855

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96 (Leve Item ID 19977)
$ws.Cells.Item(96, 8).Value = 1380  # H96: 3000 -> 1380
$ws.Cells.Item(96, 9).Value = 1250  # I96: 3000 -> 1250
$ws.Cells.Item(96, 10).Value = 1466.6666  # J96: 0 -> 1466.6666
$ws.Cells.Item(96, 11).Value = 1250  # K96: 3000 -> 1250
$ws.Cells.Item(96, 12).Value = 1466.6666  # L96: 0 -> 1466.6666
$ws.Cells.Item(96, 13).Value = 123  # M96: -1627 -> 123
$ws.Cells.Item(96, 14).Value = -4212.6666  # N96: None -> -4212.6666
# Row 123 (Leve Item ID 34127)
$ws.Cells.Item(123, 8).Value = 30261.125  # H123: 24427.77 -> 30261.125
$ws.Cells.Item(123, 10).Value = 30261.125  # J123: 24427.77 -> 30261.125
$ws.Cells.Item(123, 12).Value = 30261.125  # L123: 24427.77 -> 30261.125
$ws.Cells.Item(123, 14).Value = -40061.125  # N123: -34227.77 -> -40061.125
# Row 125 (Leve Item ID 34276)
$ws.Cells.Item(125, 8).Value = 88857.22  # H125: 86043 -> 88857.22
$ws.Cells.Item(125, 10).Value = 88857.22  # J125: 86043 -> 88857.22
$ws.Cells.Item(125, 12).Value = 88857.22  # L125: 86043 -> 88857.22
$ws.Cells.Item(125, 14).Value = -98697.22  # N125: -95883 -> -98697.22
